$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1134.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1134.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 1134.5
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = -1946.5

$ws.Range("H91").Value = 1134.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1134.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 1134.5
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = -3942.5

$ws.Range("H92").Value = 1379467.9
$ws.Range("I92").Value = 781818.3
$ws.Range("J92").Value = 2233253
$ws.Range("K92").Value = 781818.3
$ws.Range("L92").Value = 2233253
$ws.Range("M92").Value = -780570.3
$ws.Range("N92").Value = -2235749

$ws.Range("H113").Value = 4143.364
$ws.Range("I113").Value = 2763.1667
$ws.Range("J113").Value = 5799.6
$ws.Range("K113").Value = 2763.1667
$ws.Range("L113").Value = 5799.6
$ws.Range("M113").Value = 490.8332999999998
$ws.Range("N113").Value = -12307.6

$ws.Range("H116").Value = 7301.0625
$ws.Range("I116").Value = 6785.2856
$ws.Range("J116").Value = 7702.222
$ws.Range("K116").Value = 6785.2856
$ws.Range("L116").Value = 7702.222
$ws.Range("M116").Value = -3343.2856
$ws.Range("N116").Value = -14586.222

$ws.Range("H132").Value = 256124.75
$ws.Range("I132").Value = 256124.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 768374.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -765844.25
$ws.Range("N132").Value = $null

$ws.Range("H137").Value = 7695362
$ws.Range("I137").Value = 1937.3334
$ws.Range("J137").Value = 14289726
$ws.Range("K137").Value = 5812.0002
$ws.Range("L137").Value = 42869178
$ws.Range("M137").Value = -3262.0002
$ws.Range("N137").Value = -42874278

$ws.Range("H138").Value = 5671.457
$ws.Range("I138").Value = 8575.869000000001
$ws.Range("J138").Value = 4250.149
$ws.Range("K138").Value = 25727.607
$ws.Range("L138").Value = 12750.447
$ws.Range("M138").Value = -20587.607
$ws.Range("N138").Value = -23030.447

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12
$ws.Range("I32").Value = 12
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 12
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 275
$ws.Range("N32").Value = $null

$ws.Range("H61").Value = 1054704.6
$ws.Range("I61").Value = 3737.8044
$ws.Range("J61").Value = 5889152.5
$ws.Range("K61").Value = 3737.8044
$ws.Range("L61").Value = 5889152.5
$ws.Range("M61").Value = -3525.8044
$ws.Range("N61").Value = -5889576.5

$ws.Range("H74").Value = 657344.8
$ws.Range("I74").Value = 1519.8667
$ws.Range("J74").Value = 1175101.4
$ws.Range("K74").Value = 1519.8667
$ws.Range("L74").Value = 1175101.4
$ws.Range("M74").Value = -645.8667
$ws.Range("N74").Value = -1176849.4

$ws.Range("H77").Value = 657344.8
$ws.Range("I77").Value = 1519.8667
$ws.Range("J77").Value = 1175101.4
$ws.Range("K77").Value = 7599.333500000001
$ws.Range("L77").Value = 5875507
$ws.Range("M77").Value = -3231.333500000001
$ws.Range("N77").Value = -5884243

$ws.Range("H97").Value = 12631.556
$ws.Range("I97").Value = 14023
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 14023
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -13527
$ws.Range("N97").Value = -2492

$ws.Range("H136").Value = 1054704.6
$ws.Range("I136").Value = 3737.8044
$ws.Range("J136").Value = 5889152.5
$ws.Range("K136").Value = 11213.4132
$ws.Range("L136").Value = 17667457.5
$ws.Range("M136").Value = -8663.413199999999
$ws.Range("N136").Value = -17672557.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1153.8948
$ws.Range("I20").Value = 1086.5
$ws.Range("J20").Value = 1185
$ws.Range("K20").Value = 1086.5
$ws.Range("L20").Value = 1185
$ws.Range("M20").Value = -839.5
$ws.Range("N20").Value = -1679

$ws.Range("H86").Value = 1723.6333
$ws.Range("I86").Value = 1036.84
$ws.Range("J86").Value = 5157.6
$ws.Range("K86").Value = 1036.84
$ws.Range("L86").Value = 5157.6
$ws.Range("M86").Value = 86.16000000000008
$ws.Range("N86").Value = -7403.6

$ws.Range("H89").Value = 1723.6333
$ws.Range("I89").Value = 1036.84
$ws.Range("J89").Value = 5157.6
$ws.Range("K89").Value = 5184.2
$ws.Range("L89").Value = 25788
$ws.Range("M89").Value = 431.8000000000002
$ws.Range("N89").Value = -37020

$ws.Range("H99").Value = 8368.85
$ws.Range("I99").Value = 15408.667
$ws.Range("J99").Value = 2609
$ws.Range("K99").Value = 15408.667
$ws.Range("L99").Value = 2609
$ws.Range("M99").Value = -13910.667
$ws.Range("N99").Value = -5605

$ws.Range("H134").Value = 20001722
$ws.Range("I134").Value = 1234.5834
$ws.Range("J134").Value = 100003670
$ws.Range("K134").Value = 3703.7502
$ws.Range("L134").Value = 300011010
$ws.Range("M134").Value = -1168.7502
$ws.Range("N134").Value = -300016080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4131.933
$ws.Range("I16").Value = 2666.2222
$ws.Range("J16").Value = 6330.5
$ws.Range("K16").Value = 2666.2222
$ws.Range("L16").Value = 6330.5
$ws.Range("M16").Value = -2379.2222
$ws.Range("N16").Value = -6904.5

$ws.Range("H86").Value = 25997
$ws.Range("I86").Value = 11330.667
$ws.Range("J86").Value = 40663.332
$ws.Range("K86").Value = 11330.667
$ws.Range("L86").Value = 40663.332
$ws.Range("M86").Value = -10207.667
$ws.Range("N86").Value = -42909.332

$ws.Range("H89").Value = 25997
$ws.Range("I89").Value = 11330.667
$ws.Range("J89").Value = 40663.332
$ws.Range("K89").Value = 56653.335
$ws.Range("L89").Value = 203316.66
$ws.Range("M89").Value = -51037.335
$ws.Range("N89").Value = -214548.66

$ws.Range("H105").Value = 3596.2307
$ws.Range("I105").Value = 2972.3333
$ws.Range("K105").Value = 2972.3333
$ws.Range("M105").Value = -1225.3333

$ws.Range("H113").Value = 4131.933
$ws.Range("I113").Value = 2666.2222
$ws.Range("J113").Value = 6330.5
$ws.Range("K113").Value = 2666.2222
$ws.Range("L113").Value = 6330.5
$ws.Range("M113").Value = -496.2222000000002
$ws.Range("N113").Value = -10670.5

$ws.Range("H122").Value = 3028.8667
$ws.Range("I122").Value = 3398.625
$ws.Range("J122").Value = 2606.2856
$ws.Range("K122").Value = 10195.875
$ws.Range("L122").Value = 7818.8568
$ws.Range("M122").Value = -7745.875
$ws.Range("N122").Value = -12718.8568

$ws.Range("H132").Value = 2712.6
$ws.Range("I132").Value = 2550.5898
$ws.Range("J132").Value = 3765.6667
$ws.Range("K132").Value = 7651.769400000001
$ws.Range("L132").Value = 11297.0001
$ws.Range("M132").Value = -5121.769400000001
$ws.Range("N132").Value = -16357.0001

$ws.Range("H134").Value = 2241.7368
$ws.Range("I134").Value = 1746
$ws.Range("J134").Value = 4100.75
$ws.Range("K134").Value = 5238
$ws.Range("L134").Value = 12302.25
$ws.Range("M134").Value = -2703
$ws.Range("N134").Value = -17372.25

$ws.Range("H141").Value = 564499.1
$ws.Range("I141").Value = 949999
$ws.Range("J141").Value = 509427.72
$ws.Range("K141").Value = 949999
$ws.Range("L141").Value = 509427.72
$ws.Range("M141").Value = -944819
$ws.Range("N141").Value = -519787.72

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8749555
$ws.Range("I80").Value = 167111.31
$ws.Range("J80").Value = 21809794
$ws.Range("K80").Value = 167111.31
$ws.Range("L80").Value = 21809794
$ws.Range("M80").Value = -166113.31
$ws.Range("N80").Value = -21811790

$ws.Range("H83").Value = 8749555
$ws.Range("I83").Value = 167111.31
$ws.Range("J83").Value = 21809794
$ws.Range("K83").Value = 835556.55
$ws.Range("L83").Value = 109048970
$ws.Range("M83").Value = -830564.55
$ws.Range("N83").Value = -109058954

$ws.Range("H132").Value = 560417.25
$ws.Range("I132").Value = 2555.4
$ws.Range("J132").Value = 2335432.2
$ws.Range("K132").Value = 7666.200000000001
$ws.Range("L132").Value = 7006296.600000001
$ws.Range("M132").Value = -5136.200000000001
$ws.Range("N132").Value = -7011356.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1686260
$ws.Range("I40").Value = 2059877.8
$ws.Range("J40").Value = 4979.6665
$ws.Range("K40").Value = 2059877.8
$ws.Range("L40").Value = 4979.6665
$ws.Range("M40").Value = -2059741.8
$ws.Range("N40").Value = -5251.6665

$ws.Range("H93").Value = 3656.3076
$ws.Range("I93").Value = 3412
$ws.Range("J93").Value = 5000
$ws.Range("K93").Value = 3412
$ws.Range("L93").Value = 5000
$ws.Range("M93").Value = -2164
$ws.Range("N93").Value = -7496

$ws.Range("H100").Value = 4402.8335
$ws.Range("I100").Value = 2854.25
$ws.Range("K100").Value = 2854.25
$ws.Range("M100").Value = -2313.25

$ws.Range("H132").Value = 13947.077
$ws.Range("I132").Value = 5199.6
$ws.Range("J132").Value = 19414.25
$ws.Range("K132").Value = 15598.8
$ws.Range("L132").Value = 58242.75
$ws.Range("M132").Value = -13068.8
$ws.Range("N132").Value = -63302.75

$ws.Range("H136").Value = 4073.7693
$ws.Range("I136").Value = 2342.7778
$ws.Range("J136").Value = 7968.5
$ws.Range("K136").Value = 7028.3334
$ws.Range("L136").Value = 23905.5
$ws.Range("M136").Value = -4478.3334
$ws.Range("N136").Value = -29005.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1787.7059
$ws.Range("I122").Value = 1787.7059
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5363.1177
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2913.1177
$ws.Range("N122").Value = $null

$ws.Range("H132").Value = 28409.264
$ws.Range("I132").Value = 34708.4
$ws.Range("J132").Value = 4787.5
$ws.Range("K132").Value = 104125.2
$ws.Range("L132").Value = 14362.5
$ws.Range("M132").Value = -101595.2
$ws.Range("N132").Value = -19422.5
